$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 467415
$ws.Range("J17").Value = 467415
$ws.Range("L17").Value = 1402245
$ws.Range("N17").Value = -1402581

$ws.Range("H132").Value = 1819979.4
$ws.Range("I132").Value = 1680.5814
$ws.Range("J132").Value = 8335550
$ws.Range("K132").Value = 5041.7442
$ws.Range("L132").Value = 25006650
$ws.Range("M132").Value = -2511.7442
$ws.Range("N132").Value = -25011710

$ws.Range("H135").Value = 1462.0834
$ws.Range("I135").Value = 1499.5238
$ws.Range("J135").Value = 1200
$ws.Range("K135").Value = 13495.7142
$ws.Range("L135").Value = 10800
$ws.Range("M135").Value = -10960.7142
$ws.Range("N135").Value = -15870

$ws.Range("H138").Value = 3859.3298
$ws.Range("I138").Value = 2317.682
$ws.Range("J138").Value = 4330.3887
$ws.Range("K138").Value = 6953.045999999999
$ws.Range("L138").Value = 12991.1661
$ws.Range("M138").Value = -1813.045999999999
$ws.Range("N138").Value = -23271.1661

$ws.Range("H141").Value = 749.375
$ws.Range("I141").Value = 618
$ws.Range("J141").Value = 968.3333
$ws.Range("K141").Value = 1854
$ws.Range("L141").Value = 2904.9999
$ws.Range("M141").Value = 3326
$ws.Range("N141").Value = -13264.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 116.666664
$ws.Range("I5").Value = 80
$ws.Range("J5").Value = 300
$ws.Range("K5").Value = 80
$ws.Range("L5").Value = 300
$ws.Range("M5").Value = 32
$ws.Range("N5").Value = -524

$ws.Range("H11").Value = 1000
$ws.Range("I11").Value = 1000
$ws.Range("K11").Value = 1000
$ws.Range("M11").Value = -856

$ws.Range("H61").Value = 1630.1428
$ws.Range("I61").Value = 1670.9231
$ws.Range("J61").Value = 1100
$ws.Range("K61").Value = 1670.9231
$ws.Range("L61").Value = 1100
$ws.Range("M61").Value = -1458.9231
$ws.Range("N61").Value = -1524

$ws.Range("H136").Value = 1630.1428
$ws.Range("I136").Value = 1670.9231
$ws.Range("J136").Value = 1100
$ws.Range("K136").Value = 5012.7693
$ws.Range("L136").Value = 3300
$ws.Range("M136").Value = -2462.7693
$ws.Range("N136").Value = -8400

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 116.666664
$ws.Range("I4").Value = 80
$ws.Range("J4").Value = 300
$ws.Range("K4").Value = 80
$ws.Range("L4").Value = 300
$ws.Range("M4").Value = 35
$ws.Range("N4").Value = -530

$ws.Range("H23").Value = 60014
$ws.Range("J23").Value = 60014
$ws.Range("L23").Value = 60014
$ws.Range("N23").Value = -60580

$ws.Range("H134").Value = 14278.077
$ws.Range("I134").Value = 1253.0615
$ws.Range("J134").Value = 79403.16
$ws.Range("K134").Value = 3759.1845
$ws.Range("L134").Value = 238209.48
$ws.Range("M134").Value = -1224.1845
$ws.Range("N134").Value = -243279.48

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 200
$ws.Range("I17").Value = 200
$ws.Range("K17").Value = 200
$ws.Range("M17").Value = -26

$ws.Range("H31").Value = 2465.4792
$ws.Range("I31").Value = 2034.3684
$ws.Range("K31").Value = 2034.3684
$ws.Range("M31").Value = -1739.3684

$ws.Range("H34").Value = 2465.4792
$ws.Range("I34").Value = 2034.3684
$ws.Range("K34").Value = 2034.3684
$ws.Range("M34").Value = -1832.3684

$ws.Range("H58").Value = 3633.027
$ws.Range("I58").Value = 800.6875
$ws.Range("J58").Value = 21760
$ws.Range("K58").Value = 800.6875
$ws.Range("L58").Value = 21760
$ws.Range("M58").Value = -597.6875
$ws.Range("N58").Value = -22166

$ws.Range("H132").Value = 1384.2142
$ws.Range("I132").Value = 937.8276
$ws.Range("J132").Value = 2380
$ws.Range("K132").Value = 2813.4828
$ws.Range("L132").Value = 7140
$ws.Range("M132").Value = -283.4827999999998
$ws.Range("N132").Value = -12200

$ws.Range("H136").Value = 3633.027
$ws.Range("I136").Value = 800.6875
$ws.Range("J136").Value = 21760
$ws.Range("K136").Value = 2402.0625
$ws.Range("L136").Value = 65280
$ws.Range("M136").Value = 147.9375
$ws.Range("N136").Value = -70380

$ws.Range("H137").Value = 78570
$ws.Range("J137").Value = 78570
$ws.Range("L137").Value = 78570
$ws.Range("N137").Value = -88770

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 160
$ws.Range("I40").Value = 160
$ws.Range("K40").Value = 640
$ws.Range("M40").Value = -571

$ws.Range("H80").Value = 8155.6665
$ws.Range("J80").Value = 8666.666999999999
$ws.Range("L80").Value = 26000.001
$ws.Range("N80").Value = -27872.001

$ws.Range("H83").Value = 8155.6665
$ws.Range("J83").Value = 8666.666999999999
$ws.Range("L83").Value = 78000.003
$ws.Range("N83").Value = -87360.003

$ws.Range("H122").Value = 28371.27
$ws.Range("I122").Value = 548.1667
$ws.Range("J122").Value = 33756.387
$ws.Range("K122").Value = 4933.5003
$ws.Range("L122").Value = 303807.483
$ws.Range("M122").Value = -2483.5003
$ws.Range("N122").Value = -308707.483

$ws.Range("H125").Value = 2553.4443
$ws.Range("I125").Value = 982.6667
$ws.Range("J125").Value = 3338.8333
$ws.Range("K125").Value = 2948.0001
$ws.Range("L125").Value = 10016.4999
$ws.Range("M125").Value = 1971.9999
$ws.Range("N125").Value = -19856.4999

$ws.Range("H131").Value = 78554.38
$ws.Range("I131").Value = 112717.78
$ws.Range("J131").Value = 60467.883
$ws.Range("K131").Value = 338153.34
$ws.Range("L131").Value = 181403.649
$ws.Range("M131").Value = -333113.34
$ws.Range("N131").Value = -191483.649

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 5786.5
$ws.Range("I126").Value = 4132
$ws.Range("K126").Value = 12396
$ws.Range("M126").Value = -9926

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2820.8
$ws.Range("I122").Value = 2620.1365
$ws.Range("K122").Value = 7860.4095
$ws.Range("M122").Value = -5410.4095

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 954.38464
$ws.Range("I126").Value = 927.9091
$ws.Range("K126").Value = 2783.7273
$ws.Range("M126").Value = -313.7273

$ws.Range("H132").Value = 5000.5864
$ws.Range("I132").Value = 1416
$ws.Range("J132").Value = 11811.3
$ws.Range("K132").Value = 4248
$ws.Range("L132").Value = 35433.89999999999
$ws.Range("M132").Value = -1718
$ws.Range("N132").Value = -40493.89999999999

$ws.Range("H136").Value = 3555.9412
$ws.Range("I136").Value = 576.6
$ws.Range("J136").Value = 25901
$ws.Range("K136").Value = 1729.8
$ws.Range("L136").Value = 77703
$ws.Range("M136").Value = 820.1999999999998
$ws.Range("N136").Value = -82803
